$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H8").Value = 585.2857
$ws.Range("I8").Value = 585.2857
$ws.Range("J8").Value = 0
$ws.Range("K8").Value = 1755.8571
$ws.Range("L8").Value = 0
$ws.Range("M8").ClearContents()
$ws.Range("N8").Value = -1616.8571
$ws.Range("H15").Value = 250682.75
$ws.Range("I15").Value = 250682.75
$ws.Range("K15").Value = 752048.25
$ws.Range("M15").Value = -751879.25
$ws.Range("H132").Value = 0
$ws.Range("I132").Value = 0
$ws.Range("K132").Value = 0
$ws.Range("M132").ClearContents()
$ws.Range("H137").Value = 1687.4348
$ws.Range("I137").Value = 1316.9474
$ws.Range("J137").Value = 3447.25
$ws.Range("K137").Value = 3950.8422
$ws.Range("L137").Value = 10341.75
$ws.Range("M137").Value = -1400.8422
$ws.Range("N137").Value = -15441.75

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 1103.6471
$ws.Range("I61").Value = 988.5
$ws.Range("K61").Value = 988.5
$ws.Range("M61").Value = -776.5
$ws.Range("H74").Value = 1261.9131
$ws.Range("I74").Value = 1399.6
$ws.Range("J74").Value = 1003.75
$ws.Range("K74").Value = 1399.6
$ws.Range("L74").Value = 1003.75
$ws.Range("M74").Value = -525.5999999999999
$ws.Range("N74").Value = -2751.75
$ws.Range("H77").Value = 1261.9131
$ws.Range("I77").Value = 1399.6
$ws.Range("J77").Value = 1003.75
$ws.Range("K77").Value = 6998
$ws.Range("L77").Value = 5018.75
$ws.Range("M77").Value = -2630
$ws.Range("N77").Value = -13754.75
$ws.Range("H122").Value = 1777.8334
$ws.Range("I122").Value = 1135
$ws.Range("J122").Value = 3063.5
$ws.Range("K122").Value = 3405
$ws.Range("L122").Value = 9190.5
$ws.Range("M122").Value = -955
$ws.Range("N122").Value = -14090.5
$ws.Range("H132").Value = 2283.7886
$ws.Range("I132").Value = 1069.7188
$ws.Range("K132").Value = 3209.1564
$ws.Range("M132").Value = -679.1564000000003
$ws.Range("H136").Value = 1103.6471
$ws.Range("I136").Value = 988.5
$ws.Range("K136").Value = 2965.5
$ws.Range("M136").Value = -415.5

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H6").Value = 23980
$ws.Range("J6").Value = 23980
$ws.Range("L6").Value = 23980
$ws.Range("N6").Value = -24206
$ws.Range("H105").Value = 2785.5881
$ws.Range("I105").Value = 2510.6428
$ws.Range("J105").Value = 4068.6667
$ws.Range("K105").Value = 2510.6428
$ws.Range("L105").Value = 4068.6667
$ws.Range("M105").Value = -763.6428000000001
$ws.Range("N105").Value = -7562.6667
$ws.Range("H134").Value = 2070.2974
$ws.Range("I134").Value = 1224.3043
$ws.Range("J134").Value = 3460.1428
$ws.Range("K134").Value = 3672.9129
$ws.Range("L134").Value = 10380.4284
$ws.Range("M134").Value = -1137.9129
$ws.Range("N134").Value = -15450.4284

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H10").Value = 297.5
$ws.Range("I10").Value = 297.5
$ws.Range("K10").Value = 297.5
$ws.Range("M10").Value = -158.5
$ws.Range("H31").Value = 2323.2744
$ws.Range("I31").Value = 1909.0625
$ws.Range("J31").Value = 3020.8948
$ws.Range("K31").Value = 1909.0625
$ws.Range("L31").Value = 3020.8948
$ws.Range("M31").Value = -1614.0625
$ws.Range("N31").Value = -3610.8948
$ws.Range("H34").Value = 2323.2744
$ws.Range("I34").Value = 1909.0625
$ws.Range("J34").Value = 3020.8948
$ws.Range("K34").Value = 1909.0625
$ws.Range("L34").Value = 3020.8948
$ws.Range("M34").Value = -1707.0625
$ws.Range("N34").Value = -3424.8948
$ws.Range("H58").Value = 2570.8965
$ws.Range("I58").Value = 1918.8
$ws.Range("J58").Value = 2914.1052
$ws.Range("K58").Value = 1918.8
$ws.Range("L58").Value = 2914.1052
$ws.Range("M58").Value = -1715.8
$ws.Range("N58").Value = -3320.1052
$ws.Range("H109").Value = 21585
$ws.Range("J109").Value = 21585
$ws.Range("L109").Value = 21585
$ws.Range("N109").Value = -23665
$ws.Range("H132").Value = 2764.3684
$ws.Range("I132").Value = 2040.3846
$ws.Range("J132").Value = 4333
$ws.Range("K132").Value = 6121.1538
$ws.Range("L132").Value = 12999
$ws.Range("M132").Value = -3591.1538
$ws.Range("N132").Value = -18059
$ws.Range("H134").Value = 2716
$ws.Range("I134").Value = 2625.7646
$ws.Range("J134").Value = 3022.8
$ws.Range("K134").Value = 7877.293799999999
$ws.Range("L134").Value = 9068.400000000001
$ws.Range("M134").Value = -5342.293799999999
$ws.Range("N134").Value = -14138.4
$ws.Range("H136").Value = 2570.8965
$ws.Range("I136").Value = 1918.8
$ws.Range("J136").Value = 2914.1052
$ws.Range("K136").Value = 5756.4
$ws.Range("L136").Value = 8742.3156
$ws.Range("M136").Value = -3206.4
$ws.Range("N136").Value = -13842.3156

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H96").Value = 10750
$ws.Range("J96").Value = 10750
$ws.Range("L96").Value = 10750
$ws.Range("N96").Value = -16242
$ws.Range("H132").Value = 2408.1052
$ws.Range("I132").Value = 1893.9166
$ws.Range("J132").Value = 3289.5715
$ws.Range("K132").Value = 5681.7498
$ws.Range("L132").Value = 9868.7145
$ws.Range("M132").Value = -3151.7498
$ws.Range("N132").Value = -14928.7145

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H122").Value = 3543.375
$ws.Range("I122").Value = 2278.2778
$ws.Range("J122").Value = 5169.9287
$ws.Range("K122").Value = 6834.8334
$ws.Range("L122").Value = 15509.7861
$ws.Range("M122").Value = -4384.8334
$ws.Range("N122").Value = -20409.7861
$ws.Range("H132").Value = 12572.76
$ws.Range("I132").Value = 17879.072
$ws.Range("J132").Value = 5819.273
$ws.Range("K132").Value = 53637.216
$ws.Range("L132").Value = 17457.819
$ws.Range("M132").Value = -51107.216
$ws.Range("N132").Value = -22517.819
$ws.Range("H136").Value = 4005.4443
$ws.Range("I136").Value = 4656.0835
$ws.Range("K136").Value = 13968.2505
$ws.Range("M136").Value = -11418.2505

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H27").Value = 32162.5
$ws.Range("J27").Value = 32162.5
$ws.Range("L27").Value = 32162.5
$ws.Range("N27").Value = -32300.5
$ws.Range("H115").Value = 40000
$ws.Range("J115").Value = 40000
$ws.Range("L115").Value = 40000
$ws.Range("N115").Value = -43134
$ws.Range("H132").Value = 2223.375
$ws.Range("I132").Value = 1843.15
$ws.Range("J132").Value = 4124.5
$ws.Range("K132").Value = 5529.450000000001
$ws.Range("L132").Value = 12373.5
$ws.Range("M132").Value = -2999.450000000001
$ws.Range("N132").Value = -17433.5
$ws.Range("H136").Value = 13892756
$ws.Range("I136").Value = 23810532
$ws.Range("K136").Value = 71431596
$ws.Range("M136").Value = -71429046
